$d = $word.ActiveDocument
Write-Host "Test: " $d.Content.Text.Length
